# "adam & angel fall 2018"
#
# The responsibility chart tracked a handful of tasks with a "status"
# placeholder (shared strings "Waiting on" / "In Progress", shown with the
# builtin "Bad" (red) / "Good" (green) cell styles) in place of a real
# completion date. Adam & Angel wrapped those items up, so each of those
# cells now gets its actual completion date, written with the same plain
# date formatting used by every other row in the Completion Date column.
# That leaves the "Bad"/"Good" cell styles unused, so they are removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows whose Completion Date cell showed "In Progress" (style "Good") -
# these tasks finished 12/1/2018.
$finishedInProgress = @(11, 14, 15, 16, 17, 18)
foreach ($r in $finishedInProgress) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Style = "Normal"
    $cell.Value = 43435
    $cell.NumberFormat = "d-mmm"
}

# Rows whose Completion Date cell showed "Waiting on" (style "Bad") - these
# tasks finished on the dates below.
$finishedWaiting = @{
    19 = 43421
    20 = 43422
    21 = 43423
    22 = 43424
    23 = 43425
    24 = 43423
    25 = 43423
}
foreach ($r in $finishedWaiting.Keys) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Style = "Normal"
    $cell.Value = $finishedWaiting[$r]
    $cell.NumberFormat = "d-mmm"
}

# "Bad"/"Good" cell styles are no longer referenced by any cell - drop them.
$wb.Styles.Item("Bad").Delete()
$wb.Styles.Item("Good").Delete()

# Leave the selection/scroll position where Adam last left off editing.
$ws.Activate()
$ws.Range("D20").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
